$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 43.21270066666667
$ws.Range("H2").Value = 129.638102
$ws.Range("I2").Value = 0.1487696778665633
$ws.Range("J2").Value = 0.1487696778665633
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 911.0021034881091
$ws.Range("R2").Value = 8199.018931392982
$ws.Range("S2").Value = 0.008503764265841563
$ws.Range("T2").Value = 0.008503764265841561

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.21270066666667
$ws.Range("H3").Value = 129.638102
$ws.Range("I3").Value = 0.1487696778665633
$ws.Range("J3").Value = 0.1487696778665633
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 13032.95694535483
$ws.Range("R3").Value = 117296.6125081935
$ws.Range("S3").Value = 0.1216563530707661
$ws.Range("T3").Value = 0.1216563530707661

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 43.21270066666667
$ws.Range("H4").Value = 129.638102
$ws.Range("I4").Value = 0.1487696778665633
$ws.Range("J4").Value = 0.1487696778665633
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 1993.6287340276
$ws.Range("R4").Value = 17942.6586062484
$ws.Range("S4").Value = 0.01860956052995562
$ws.Range("T4").Value = 0.01860956052995561

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 237.7114053333333
$ws.Range("H5").Value = 713.134216
$ws.Range("I5").Value = 0.8183762794517323
$ws.Range("J5").Value = 0.8183762794517323
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 5011.387553678806
$ws.Range("R5").Value = 45102.48798310925
$ws.Range("S5").Value = 0.04677888035393898
$ws.Range("T5").Value = 0.04677888035393898

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 237.7114053333333
$ws.Range("H6").Value = 713.134216
$ws.Range("I6").Value = 0.8183762794517323
$ws.Range("J6").Value = 0.8183762794517323
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 71693.7951882956
$ws.Range("R6").Value = 645244.1566946603
$ws.Range("S6").Value = 0.669226921947222
$ws.Range("T6").Value = 0.669226921947222

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 237.7114053333333
$ws.Range("H7").Value = 713.134216
$ws.Range("I7").Value = 0.8183762794517323
$ws.Range("J7").Value = 0.8183762794517323
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 10966.87503366753
$ws.Range("R7").Value = 98701.8753030078
$ws.Range("S7").Value = 0.1023704771505714
$ws.Range("T7").Value = 0.1023704771505714

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.543019333333334
$ws.Range("H8").Value = 28.629058
$ws.Range("I8").Value = 0.03285404268170446
$ws.Range("J8").Value = 0.03285404268170446
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 201.1841554027309
$ws.Range("R8").Value = 1810.657398624578
$ws.Range("S8").Value = 0.001877956840073958
$ws.Range("T8").Value = 0.001877956840073958

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.543019333333334
$ws.Range("H9").Value = 28.629058
$ws.Range("I9").Value = 0.03285404268170446
$ws.Range("J9").Value = 0.03285404268170446
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 2878.176049662207
$ws.Range("R9").Value = 25903.58444695987
$ws.Range("S9").Value = 0.02686638213919116
$ws.Range("T9").Value = 0.02686638213919116

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.543019333333334
$ws.Range("H10").Value = 28.629058
$ws.Range("I10").Value = 0.03285404268170446
$ws.Range("J10").Value = 0.03285404268170446
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 440.2695795171603
$ws.Range("R10").Value = 3962.426215654442
$ws.Range("S10").Value = 0.00410970370243935
$ws.Range("T10").Value = 0.00410970370243935

